# Auto-update draw results: append the 2025-09-22 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Columns A, C and E hold values that look numeric/date-like ("2025-09-22",
# "250922", the ISO timestamp) but must stay as plain text, matching the
# rest of the table. A leading quote forces Excel to store them as text
# instead of coercing them into dates/numbers.
$ws.Range("A$row").Value = "'2025-09-22"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "'250922"
$ws.Range("D$row").Value = "5-7-4-2"
$ws.Range("E$row").Value = "'2025-09-22T21:37:07.517+04:00"
